# TDexcel_SkillBoard.xlsx update
# Commit: Added Tests/test_DataDrivenIteration.py, Utilities/get_ddt(),
#         get_ddt_iteration(), Config/file for ddt
#
# The only functional data change in the target workbook is the test-case
# label in row 6 of the data-driven-testing table: "tc_5" becomes
# "Sign_Up" (a new Sign-Up test case row replacing the old tc_5 entry).
# The active selection also moved to A7 before the file was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Sign_Up"

$ws.Range("A7").Select()
